# Auto-generated Excel COM-interop edit script
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# across several worksheets per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1611.6
$ws.Range("I111").Value = 1279
$ws.Range("J111").Value = 1833.3334
$ws.Range("K111").Value = 3837
$ws.Range("L111").Value = 5500.0002
$ws.Range("M111").Value = -770
$ws.Range("N111").Value = -11634.0002
$ws.Range("H115").Value = 1831.75
$ws.Range("I115").Value = 1184
$ws.Range("J115").Value = 2220.4
$ws.Range("K115").Value = 3552
$ws.Range("L115").Value = 6661.200000000001
$ws.Range("M115").Value = -1985
$ws.Range("N115").Value = -9795.200000000001
$ws.Range("H129").Value = 901.8
$ws.Range("I129").Value = 500
$ws.Range("J129").Value = 963.61536
$ws.Range("K129").Value = 1500
$ws.Range("L129").Value = 2890.84608
$ws.Range("M129").Value = 3500
$ws.Range("N129").Value = -12890.84608
$ws.Range("H133").Value = 48377.273
$ws.Range("J133").Value = 48377.273
$ws.Range("L133").Value = 48377.273
$ws.Range("N133").Value = -58497.273
$ws.Range("H138").Value = 5752.04
$ws.Range("I138").Value = 968.3077
$ws.Range("J138").Value = 7432.811
$ws.Range("K138").Value = 2904.9231
$ws.Range("L138").Value = 22298.433
$ws.Range("M138").Value = 2235.0769
$ws.Range("N138").Value = -32578.433

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1350.2727
$ws.Range("I2").Value = 1306.8334
$ws.Range("J2").Value = 1402.4
$ws.Range("K2").Value = 1306.8334
$ws.Range("L2").Value = 1402.4
$ws.Range("M2").Value = -1193.8334
$ws.Range("N2").Value = -1628.4
$ws.Range("H45").Value = 1641.2222
$ws.Range("I45").Value = 1726.4
$ws.Range("J45").Value = 1534.75
$ws.Range("K45").Value = 1726.4
$ws.Range("L45").Value = 1534.75
$ws.Range("M45").Value = -1349.4
$ws.Range("N45").Value = -2288.75
$ws.Range("H60").Value = 10051
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H116").Value = 1350.2727
$ws.Range("I116").Value = 1306.8334
$ws.Range("J116").Value = 1402.4
$ws.Range("K116").Value = 1306.8334
$ws.Range("L116").Value = 1402.4
$ws.Range("M116").Value = 987.1666
$ws.Range("N116").Value = -5990.4
$ws.Range("H122").Value = 1803.931
$ws.Range("I122").Value = 1115.9375
$ws.Range("J122").Value = 2650.6924
$ws.Range("K122").Value = 3347.8125
$ws.Range("L122").Value = 7952.0772
$ws.Range("M122").Value = -897.8125
$ws.Range("N122").Value = -12852.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1350.2727
$ws.Range("I3").Value = 1306.8334
$ws.Range("J3").Value = 1402.4
$ws.Range("K3").Value = 1306.8334
$ws.Range("L3").Value = 1402.4
$ws.Range("M3").Value = -1192.8334
$ws.Range("N3").Value = -1630.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9528438
$ws.Range("I99").Value = 22225024
$ws.Range("J99").Value = 5999.1665
$ws.Range("K99").Value = 22225024
$ws.Range("L99").Value = 5999.1665
$ws.Range("M99").Value = -22223526
$ws.Range("N99").Value = -8995.166499999999
$ws.Range("H122").Value = 3990.7144
$ws.Range("I122").Value = 2147
$ws.Range("J122").Value = 8600
$ws.Range("K122").Value = 6441
$ws.Range("L122").Value = 25800
$ws.Range("M122").Value = -3991
$ws.Range("N122").Value = -30700
$ws.Range("H126").Value = 9528438
$ws.Range("I126").Value = 22225024
$ws.Range("J126").Value = 5999.1665
$ws.Range("K126").Value = 66675072
$ws.Range("L126").Value = 17997.4995
$ws.Range("M126").Value = -66672602
$ws.Range("N126").Value = -22937.4995
$ws.Range("H134").Value = 3595.3555
$ws.Range("I134").Value = 3788.8708
$ws.Range("J134").Value = 3166.8572
$ws.Range("K134").Value = 11366.6124
$ws.Range("L134").Value = 9500.571599999999
$ws.Range("M134").Value = -8831.6124
$ws.Range("N134").Value = -14570.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6263.8066
$ws.Range("I70").Value = 5817.3477
$ws.Range("J70").Value = 7547.375
$ws.Range("K70").Value = 5817.3477
$ws.Range("L70").Value = 7547.375
$ws.Range("M70").Value = -5547.3477
$ws.Range("N70").Value = -8087.375
$ws.Range("H73").Value = 6263.8066
$ws.Range("I73").Value = 5817.3477
$ws.Range("J73").Value = 7547.375
$ws.Range("K73").Value = 5817.3477
$ws.Range("L73").Value = 7547.375
$ws.Range("M73").Value = -4881.3477
$ws.Range("N73").Value = -9419.375
$ws.Range("H102").Value = 2062.9092
$ws.Range("I102").Value = 1419.2
$ws.Range("K102").Value = 1419.2
$ws.Range("M102").Value = 202.8
$ws.Range("H122").Value = 2028.25
$ws.Range("I122").Value = 1449
$ws.Range("J122").Value = 8400
$ws.Range("K122").Value = 4347
$ws.Range("L122").Value = 25200
$ws.Range("M122").Value = -1897
$ws.Range("N122").Value = -30100
$ws.Range("H126").Value = 1873.44
$ws.Range("I126").Value = 1887.9479
$ws.Range("J126").Value = 1525.25
$ws.Range("K126").Value = 5663.843699999999
$ws.Range("L126").Value = 4575.75
$ws.Range("M126").Value = -3193.843699999999
$ws.Range("N126").Value = -9515.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 93627.45
$ws.Range("I22").Value = 251625.25
$ws.Range("K22").Value = 251625.25
$ws.Range("M22").Value = -251330.25
$ws.Range("H27").Value = 93627.45
$ws.Range("I27").Value = 251625.25
$ws.Range("K27").Value = 251625.25
$ws.Range("M27").Value = -251518.25
$ws.Range("H40").Value = 6269.0415
$ws.Range("I40").Value = 4461.4
$ws.Range("K40").Value = 4461.4
$ws.Range("M40").Value = -4325.4
$ws.Range("H46").Value = 2338.0715
$ws.Range("I46").Value = 2875.25
$ws.Range("J46").Value = 2123.2
$ws.Range("K46").Value = 2875.25
$ws.Range("L46").Value = 2123.2
$ws.Range("M46").Value = -2687.25
$ws.Range("N46").Value = -2499.2
$ws.Range("H122").Value = 4227.125
$ws.Range("I122").Value = 2755.1428
$ws.Range("J122").Value = 6287.9
$ws.Range("K122").Value = 8265.428400000001
$ws.Range("L122").Value = 18863.7
$ws.Range("M122").Value = -5815.428400000001
$ws.Range("N122").Value = -23763.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H122").Value = 3108.56
$ws.Range("I122").Value = 1942.7222
$ws.Range("J122").Value = 6106.4287
$ws.Range("K122").Value = 5828.1666
$ws.Range("L122").Value = 18319.2861
$ws.Range("M122").Value = -3378.1666
$ws.Range("N122").Value = -23219.2861
$ws.Range("H126").Value = 2606.4167
$ws.Range("I126").Value = 1664.625
$ws.Range("J126").Value = 4490
$ws.Range("K126").Value = 4993.875
$ws.Range("L126").Value = 13470
$ws.Range("M126").Value = -2523.875
$ws.Range("N126").Value = -18410
